# "Separate problems from solutions"
#
# The deck originally held 4 slides:
#   1. "Problem"  (Title Slide layout, title="Problem", empty subtitle)
#   2. the actual problem diagram (Title Only layout, title placeholder
#      empty, puzzle shapes)
#   3. "Solution" (Title Slide layout, title="Solution", empty subtitle)
#   4. the solution write-up (Title and Content layout)
#
# The solution portion (slides 1, 3 and 4) is split out of this file,
# leaving only the problem diagram slide behind. Its previously-empty
# title placeholder now reads "Problem" so the slide stands on its own.

$p = $ppt.ActivePresentation

# Drop the trailing "solution" slide (step-by-step write up).
$p.Slides.Item(4).Delete()

# Drop the "Solution" title/section slide.
$p.Slides.Item(3).Delete()

# Drop the leading "Problem" title/section slide - its heading moves
# onto the remaining diagram slide below.
$p.Slides.Item(1).Delete()

# Only the problem-diagram slide remains; give it the "Problem" title.
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Problem"
